# Update "想去人数" (interested count) values on the 展览 and 全部类型 sheets
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 161
    $ws.Range("F7").Value = 659
}
